# Datos.xlsx — update AFIP "Comprobantes en Linea" login data row:
#  - CUIT (A2) refreshed
#  - Contraseña (B2) replaced with the real password value
#  - Razon social (C2) replaced with the real company/name value
# and leave the selection where the user left off (A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 20000000001
$ws.Range("B2").Value = "fakestreet123"
$ws.Range("C2").Value = "COSME FULANITO"

$ws.Range("A3").Select()
